$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range("D2","E2").NumberFormat = "@"
$ws.Range("D2").Value = '61.587.69'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D2","E2").ClearFormats()

# Row 3: 'Ethereum'
$ws.Range("D3","E3").NumberFormat = "@"
$ws.Range("D3").Value = '3.378.87'
$ws.Range("E3").Value = '  -0.36%  '
$ws.Range("D3","E3").ClearFormats()

# Row 4: 'TetherUSD'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E4").ClearFormats()

# Row 5: 'BNB'
$ws.Range("D5","E5").NumberFormat = "@"
$ws.Range("D5").Value = '576.81'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D5","E5").ClearFormats()

# Row 6: 'Solana'
$ws.Range("D6","E6").NumberFormat = "@"
$ws.Range("D6").Value = '136.39'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D6","E6").ClearFormats()

# Row 7: 'USDC'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E7").ClearFormats()

# Row 8: 'LidoStakedEther'
$ws.Range("D8","E8").NumberFormat = "@"
$ws.Range("D8").Value = '3.377.43'
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("D8","E8").ClearFormats()

# Row 9: 'XRP'
$ws.Range("D9","E9").NumberFormat = "@"
$ws.Range("D9").Value = '0.474'
$ws.Range("E9").Value = '  -0.93%  '
$ws.Range("D9","E9").ClearFormats()

# Row 10: 'Toncoin'
$ws.Range("D10","E10").NumberFormat = "@"
$ws.Range("D10").Value = '7.45'
$ws.Range("E10").Value = '  -1.77%  '
$ws.Range("D10","E10").ClearFormats()

# Row 11: 'Dogecoin'
$ws.Range("D11","E11").NumberFormat = "@"
$ws.Range("D11").Value = '0.125'
$ws.Range("E11").Value = '  +1.29%  '
$ws.Range("D11","E11").ClearFormats()

# Row 12: 'Cardano'
$ws.Range("D12","E12").NumberFormat = "@"
$ws.Range("D12").Value = '0.389'
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D12","E12").ClearFormats()

# Row 13: 'WrappedliquidstakedEther2.0'
$ws.Range("D13","E13").NumberFormat = "@"
$ws.Range("D13").Value = '3.960.78'
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D13","E13").ClearFormats()

# Row 14: 'TRON'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.75%  '
$ws.Range("E14").ClearFormats()

# Row 15: 'ShibaInu'
$ws.Range("D15","E15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000176'
$ws.Range("E15").Value = '  +0.50%  '
$ws.Range("D15","E15").ClearFormats()

# Row 16: 'WrappedEther'
$ws.Range("D16","E16").NumberFormat = "@"
$ws.Range("D16").Value = '3.386.15'
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D16","E16").ClearFormats()

# Row 17: 'Avalanche'
$ws.Range("D17","E17").NumberFormat = "@"
$ws.Range("D17").Value = '25.56'
$ws.Range("E17").Value = '  +1.41%  '
$ws.Range("D17","E17").ClearFormats()

# Row 18: 'WrappedBTC'
$ws.Range("D18","E18").NumberFormat = "@"
$ws.Range("D18").Value = '61.704.39'
$ws.Range("E18").Value = '  +0.40%  '
$ws.Range("D18","E18").ClearFormats()

# Row 19: 'Chainlink'
$ws.Range("D19","E19").NumberFormat = "@"
$ws.Range("D19").Value = '14.15'
$ws.Range("E19").Value = '  +0.54%  '
$ws.Range("D19","E19").ClearFormats()

# Row 20: 'Uniswap'
$ws.Range("D20","E20").NumberFormat = "@"
$ws.Range("D20").Value = '9.49'
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("D20","E20").ClearFormats()

# Row 21: 'Polkadot'
$ws.Range("D21","E21").NumberFormat = "@"
$ws.Range("D21").Value = '5.79'
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D21","E21").ClearFormats()

# Row 22: 'BitcoinCash'
$ws.Range("D22","E22").NumberFormat = "@"
$ws.Range("D22").Value = '378.71'
$ws.Range("E22").Value = '  +0.44%  '
$ws.Range("D22","E22").ClearFormats()

# Row 23: 'Polygon'
$ws.Range("D23","E23").NumberFormat = "@"
$ws.Range("D23").Value = '0.559'
$ws.Range("E23").Value = '  -1.83%  '
$ws.Range("D23","E23").ClearFormats()

# Row 24: 'WrappedeETH'
$ws.Range("D24","E24").NumberFormat = "@"
$ws.Range("D24").Value = '3.527.91'
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D24","E24").ClearFormats()

# Row 25: 'Dai'
$ws.Range("D25","E25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D25","E25").ClearFormats()

# Row 26: 'PEPE' -> 'Litecoin'
$ws.Range("D26","E26").NumberFormat = "@"
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = '71.09'
$ws.Range("E26").Value = '  +0.46%  '
$ws.Range("D26","E26").ClearFormats()

# Row 27: 'Litecoin' -> 'PEPE'
$ws.Range("D27","E27").NumberFormat = "@"
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").Value = '0.0000124'
$ws.Range("E27").Value = '  +4.60%  '
$ws.Range("D27","E27").ClearFormats()

# Row 28: 'Fetch.AI'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.11%  '
$ws.Range("E28").ClearFormats()

# Row 29: 'RenderToken'
$ws.Range("D29","E29").NumberFormat = "@"
$ws.Range("D29").Value = '7.57'
$ws.Range("E29").Value = '  -2.57%  '
$ws.Range("D29","E29").ClearFormats()

# Row 30: 'Binance-PegBSC-USD'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("E30").ClearFormats()

# Row 31: 'InternetComputer(DFINITY)'
$ws.Range("D31","E31").NumberFormat = "@"
$ws.Range("D31").Value = '8.16'
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("D31","E31").ClearFormats()

# Row 32: 'Kaspa'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.86%  '
$ws.Range("E32").ClearFormats()

# Row 33: 'PancakeSwap'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("E33").ClearFormats()

# Row 34: 'USDe'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("E34").ClearFormats()

# Row 35: 'EthereumClassic'
$ws.Range("D35","E35").NumberFormat = "@"
$ws.Range("D35").Value = '23.34'
$ws.Range("E35").Value = '  -0.47%  '
$ws.Range("D35","E35").ClearFormats()

# Row 36: 'NEARProtocol'
$ws.Range("D36","E36").NumberFormat = "@"
$ws.Range("D36").Value = '5.34'
$ws.Range("E36").Value = '  -4.71%  '
$ws.Range("D36","E36").ClearFormats()

# Row 37: 'ImmutableX'
$ws.Range("D37","E37").NumberFormat = "@"
$ws.Range("D37").Value = '1.55'
$ws.Range("E37").Value = '  -1.56%  '
$ws.Range("D37","E37").ClearFormats()

# Row 38: 'Aptos'
$ws.Range("D38","E38").NumberFormat = "@"
$ws.Range("D38").Value = '6.83'
$ws.Range("E38").Value = '  -1.84%  '
$ws.Range("D38","E38").ClearFormats()

# Row 39: 'Monero'
$ws.Range("D39","E39").NumberFormat = "@"
$ws.Range("D39").Value = '164.40'
$ws.Range("E39").Value = '  +0.85%  '
$ws.Range("D39","E39").ClearFormats()

# Row 40: 'Hedera'
$ws.Range("D40","E40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0784'
$ws.Range("E40").Value = '  -0.78%  '
$ws.Range("D40","E40").ClearFormats()

# Row 41: 'Mantle' -> 'FirstDigitalUSD'
$ws.Range("D41","E41").NumberFormat = "@"
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D41","E41").ClearFormats()

# Row 42: 'ONDO' -> 'Mantle'
$ws.Range("D42","E42").NumberFormat = "@"
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Value = '0.780'
$ws.Range("E42").Value = '  +2.30%  '
$ws.Range("D42","E42").ClearFormats()

# Row 43: 'FirstDigitalUSD' -> 'ONDO'
$ws.Range("D43","E43").NumberFormat = "@"
$ws.Range("B43").Value = 'ONDO'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D43").Value = '1.23'
$ws.Range("E43").Value = '  +1.16%  '
$ws.Range("D43","E43").ClearFormats()

# Row 44: 'Stacks'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +6.20%  '
$ws.Range("E44").ClearFormats()

# Row 45: 'Filecoin' -> 'EnergySwap'
$ws.Range("D45","E45").NumberFormat = "@"
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '24.86'
$ws.Range("E45").Value = '  +5.75%  '
$ws.Range("D45","E45").ClearFormats()

# Row 46: 'EnergySwap' -> 'Filecoin'
$ws.Range("D46","E46").NumberFormat = "@"
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").Value = '4.40'
$ws.Range("E46").Value = '  -0.68%  '
$ws.Range("D46","E46").ClearFormats()

# Row 47: 'OKB'
$ws.Range("D47","E47").NumberFormat = "@"
$ws.Range("D47").Value = '41.23'
$ws.Range("E47").Value = '  -0.41%  '
$ws.Range("D47","E47").ClearFormats()

# Row 48: 'Cosmos'
$ws.Range("D48","E48").NumberFormat = "@"
$ws.Range("D48").Value = '6.85'
$ws.Range("E48").Value = '  -2.00%  '
$ws.Range("D48","E48").ClearFormats()

# Row 49: 'InjectiveProtocol'
$ws.Range("D49","E49").NumberFormat = "@"
$ws.Range("D49").Value = '22.78'
$ws.Range("E49").Value = '  -1.43%  '
$ws.Range("D49","E49").ClearFormats()

# Row 50: 'Maker'
$ws.Range("D50","E50").NumberFormat = "@"
$ws.Range("D50").Value = '2.330.47'
$ws.Range("E50").Value = '  +5.53%  '
$ws.Range("D50","E50").ClearFormats()

# Row 51: 'VeChain'
$ws.Range("D51","E51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0261'
$ws.Range("E51").Value = '  +1.21%  '
$ws.Range("D51","E51").ClearFormats()
